$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gefilterd")
$ws.Activate()

# A2: add the e-mail hyperlink (display text must match the existing
# shared-string "robinkep@gmail.com", keeping the pre-existing Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:robinkep@gmail.com") | Out-Null
$ws.Range("A2").Value = "robinkep@gmail.com"
$ws.Range("A2").Style = "Hyperlink"

# D2: update the birthdate value
$ws.Range("D2").Value = 40285

# Update the visible selection
$ws.Range("A2:XFD24").Select()
